$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. The paragraph that used to read "Label the default markers with the
#    marker labels in the marker set you create in Step 3." now holds the
#    new "merging traces" sentence (the bookmark that used to sit at the end
#    of the MATLAB paragraph also moves here).
# ---------------------------------------------------------------------------
$pLabel = $d.Paragraphs.Item(7)
$pLabel.Range.Find.Execute(
    "Label the default markers with the marker labels in the marker set you create in Step 3.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "If Step 5 yields more traces than the actual number of markers, then we need to merge the traces for the same marker together manually.",
    2
)

# ---------------------------------------------------------------------------
# 2. The original "Label the default markers..." sentence becomes its own,
#    new bullet immediately after the merged paragraph (three runs collapse
#    into a single run in the new paragraph).
# ---------------------------------------------------------------------------
$pMerge = $d.Paragraphs.Item(7)
$pMerge.Range.InsertParagraphAfter()
$pNewLabel = $d.Paragraphs.Item(8)
$pNewLabel.Range.Text = "Label the default markers with the marker labels in the marker set you create in Step 3."

# ---------------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark from the end of the MATLAB paragraph to the
#    end of the merged paragraph (it stays collapsed, right before the
#    paragraph mark). Placing a bookmark exactly at a paragraph's last
#    character offset trips a position-resolution issue in this runtime, so
#    a one-character sentinel is appended first to shift that offset away
#    from the paragraph boundary, then removed once the bookmark is anchored.
# ---------------------------------------------------------------------------
$pMerge = $d.Paragraphs.Item(7)
$sentinelPos = $pMerge.Range.End - 1
$d.Range($sentinelPos, $sentinelPos).InsertAfter("X")

$pMerge = $d.Paragraphs.Item(7)
$bookmarkPos = $pMerge.Range.End - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)

$d.Bookmarks.Item("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$pMerge = $d.Paragraphs.Item(7)
$sentinelPos2 = $pMerge.Range.End - 2
$d.Range($sentinelPos2, $sentinelPos2 + 1).Delete()
